# JS-SPA-Self-Evaluation-Protocol.xlsx
# Fill in scores for two previously-blank evaluation rows and
# move the active selection/viewport to reflect where the editor
# was working next.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# "Show Buttons" row -> score 5
$ws.Range("C23").Value = 5

# "Logout" row -> score 3
$ws.Range("C29").Value = 3

# Total Score (C51 = SUM(C6:C50)) recalculates automatically: 131 -> 139

# Update the window scroll position / active selection to match
# where the author ended up after entering the new scores.
$win = $excel.ActiveWindow
$win.ScrollRow = 28
$win.ScrollColumn = 1
$ws.Range("C24").Select()
